# Apply "6 hours by turn fix": re-time the class schedule, shifting
# lunch/break rows and appending additional time slots to the end of the day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time column (A2:A17)
$times = @(
    "7:00",
    "7:50",
    "8:40",
    "9:30",
    "9:50",
    "10:40",
    "11:30",
    "12:20",
    "13:00",
    "13:50",
    "14:40",
    "15:30",
    "15:50",
    "16:40",
    "17:30",
    "18:20"
)

# New label for columns B:F, matched row-by-row with $times (row 2..17)
$labels = @(
    "-",          # 7:00  (row 2)
    "-",          # 7:50  (row 3)
    "-",          # 8:40  (row 4)
    "Intervalo",  # 9:30  (row 5)
    "-",          # 9:50  (row 6)
    "-",          # 10:40 (row 7)
    "-",          # 11:30 (row 8)
    "Almoço",     # 12:20 (row 9)
    "-",          # 13:00 (row 10)
    "-",          # 13:50 (row 11)
    "-",          # 14:40 (row 12)
    "Intervalo",  # 15:30 (row 13)
    "-",          # 15:50 (row 14)
    "-",          # 16:40 (row 15)
    "-",          # 17:30 (row 16)
    $null         # 18:20 (row 17) -- left blank
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $times[$i]

    $label = $labels[$i]
    for ($col = 2; $col -le 6; $col++) {
        if ($null -eq $label) {
            $ws.Cells.Item($row, $col).Value = ""
        } else {
            $ws.Cells.Item($row, $col).Value = $label
        }
    }
}
